$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C5:C7").Interior.ThemeColor = 5

# --- Add a "Legend" section below the chart ---
$ws.Range("B22").Value = "Legend"
$ws.Range("B22").Font.Bold = $true

# Planned swatch keeps the original yellow fill
$ws.Range("B23").Interior.Color = 65535
$ws.Range("C23").Value = "Planned"

# Actual swatch uses the new blue theme accent fill
$ws.Range("B24").Interior.ThemeColor = 5
$ws.Range("C24").Value = "Actual"
